$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC2").Value = 9.4
$ws.Range("AN2").Value = 9.800000000000001
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 2.04
$ws.Range("R2").Value = 1.4
$ws.Range("T2").Value = 1.89
$ws.Range("U2").Value = 2
$ws.Range("U5").Value = 2.02
$ws.Range("W6").Value = 3.3
$ws.Range("AB7").Value = 8.6
$ws.Range("AC7").Value = 9
$ws.Range("AG7").Value = 11
$ws.Range("I7").Value = 5.8
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 3.5
$ws.Range("O7").Value = 1.33
$ws.Range("R7").Value = 1.32
$ws.Range("S7").Value = 3.5
$ws.Range("T7").Value = 1.86
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 1.21
$ws.Range("AB9").Value = 11
$ws.Range("AC9").Value = 8
$ws.Range("F9").Value = 2.26
$ws.Range("W9").Value = 1.65
$ws.Range("F10").Value = 3.4
$ws.Range("I10").Value = 2.44
$ws.Range("J10").Value = 3.2
$ws.Range("K10").Value = 5.1
$ws.Range("R10").Value = 1.24
$ws.Range("S10").Value = 3.55
$ws.Range("I11").Value = 1.63
$ws.Range("O11").Value = 1.38
$ws.Range("R11").Value = 1.31
$ws.Range("AB12").Value = 9.4
$ws.Range("AL12").Value = 55
$ws.Range("AN12").Value = 36
$ws.Range("F12").Value = 2.8
$ws.Range("G12").Value = 2.84
$ws.Range("H12").Value = 2.92
$ws.Range("I12").Value = 2.94
$ws.Range("O12").Value = 1.44
$ws.Range("U12").Value = 1.96
$ws.Range("V12").Value = 1.51
$ws.Range("W12").Value = 1.54
$ws.Range("Y12").Value = 9.6
$ws.Range("Z12").Value = 17.5
$ws.Range("I15").Value = 4.4
$ws.Range("L15").Value = 1.31
$ws.Range("R15").Value = 1.4
$ws.Range("V15").Value = 1.3
$ws.Range("L16").Value = 1.43
$ws.Range("S16").Value = 3.7
$ws.Range("H17").Value = 5.5
$ws.Range("T17").Value = 1.96
$ws.Range("F18").Value = 3.75
$ws.Range("N18").Value = 3.25
$ws.Range("P18").Value = 1.77
$ws.Range("S18").Value = 3.9
$ws.Range("Q19").Value = 1.83
$ws.Range("N21").Value = 1.36
$ws.Range("P21").Value = 1.36
